# "fixed ion channel naming"
# On the "Gene list edited 11052021" sheet, column A cells that were tagged
# with the old ("ion Channel") label get relabeled to the corrected
# ("Ion Channel") text, for every row in the Ion-Channel gene block except
# the ones that already used the corrected label (rows 30 and 39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(31,32,33,34,35,36,37,38,40,41,42,43,44)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "ion Channel") {
        $cell.Value = "Ion Channel"
    }
}

# Update the view state: the sheet was scrolled/zoomed in further and the
# selection moved from B45 to C40.
[void]$ws.Range("C40").Select()
$excel.ActiveWindow.Zoom = 200
